# Update the "Pais" worksheet with refreshed country data and reordered
# country name rows, matching the upstream data refresh that happened
# between 04:46 and 06:03 on 24 July 2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp cell -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 06:03"

# --- Rows 31/32: Kazajistan overtakes Ecuador -------------------------------
$ws.Range("A31").Value = "Kazajistan"
$ws.Range("B31").Value = 78486
$ws.Range("C31").Value = 1687
$ws.Range("D31").Value = 48220
$ws.Range("E31").Value = 29681
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 585

$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 78148
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 33455
$ws.Range("E32").Value = 39254
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 5439

# --- Row 37: Belgica data refresh ------------------------------------------
$ws.Range("B37").Value = 64847
$ws.Range("C37").Value = 220
$ws.Range("D37").Value = 17369
$ws.Range("E37").Value = 37666
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 9812

# --- Rows 52/53: Honduras overtakes Armenia ---------------------------------
$ws.Range("A52").Value = "Honduras"
$ws.Range("B52").Value = 36902
$ws.Range("C52").Value = 800
$ws.Range("D52").Value = 4448
$ws.Range("E52").Value = 31443
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 5
$ws.Range("H52").Value = 1011

$ws.Range("A53").Value = "Armenia"
$ws.Range("B53").Value = 36162
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 25244
$ws.Range("E53").Value = 10230
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 688

# --- Row 90: Haiti data refresh ---------------------------------------------
$ws.Range("B90").Value = 7197
$ws.Range("C90").Value = 30
$ws.Range("D90").Value = 4236
$ws.Range("E90").Value = 2807
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 154

# --- Row 170: Mongolia data refresh -----------------------------------------
$ws.Range("B170").Value = 288
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 217
$ws.Range("E170").Value = 71
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0

# --- Rows 210/211: Islas Malvinas overtakes Groenlandia ---------------------
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Groenlandia"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
